# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, matching the refreshed data pull at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 325
$ws.Range("F3").Value  = 1293
$ws.Range("F5").Value  = 349
$ws.Range("F6").Value  = 3877
$ws.Range("F8").Value  = 761
$ws.Range("F9").Value  = 2286
$ws.Range("F10").Value = 341
$ws.Range("F13").Value = 178
$ws.Range("F14").Value = 180
$ws.Range("F15").Value = 2207
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 50
$ws.Range("F22").Value = 272

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 93
$ws.Range("F21").Value = 60

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2104
$ws.Range("F5").Value = 332
$ws.Range("F6").Value = 3

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 2104
$ws.Range("F5").Value  = 332
$ws.Range("F9").Value  = 325
$ws.Range("F10").Value = 1293
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 349
$ws.Range("F17").Value = 3877
$ws.Range("F22").Value = 93
$ws.Range("F23").Value = 761
$ws.Range("F24").Value = 2286
$ws.Range("F25").Value = 341
$ws.Range("F29").Value = 178
$ws.Range("F30").Value = 180
$ws.Range("F32").Value = 2207
$ws.Range("F36").Value = 17
$ws.Range("F37").Value = 50
$ws.Range("F47").Value = 60
$ws.Range("F48").Value = 272

$wb.Save()
